$wb = $excel.ActiveWorkbook

# --- Create the new "2022-Q1" sheet (moved into position after it is populated) ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q1"

# Copy header style (bold/border/center, same as other quarter sheets) from "2021-Q4"
$ref = $wb.Worksheets.Item("2021-Q4")
$ref.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$ref.Range("A2:A16").Copy($newSheet.Range("A2:A16"))

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Force text storage (matches source formatting) for size/position/ratio/value columns
$newSheet.Range("B2:G16").NumberFormat = "@"

# Fund rows
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "002408"
$newSheet.Range("C2").Value = "中信建投医改灵活配置混合A"
$newSheet.Range("D2").Value = "25.68"
$newSheet.Range("E2").Value = "94.92"
$newSheet.Range("F2").Value = "4.65"
$newSheet.Range("G2").Value = "1.1941"
$newSheet.Range("H2").Value = 8

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "000711"
$newSheet.Range("C3").Value = "嘉实医疗保健股票"
$newSheet.Range("D3").Value = "13.07"
$newSheet.Range("E3").Value = "91.70"
$newSheet.Range("F3").Value = "5.48"
$newSheet.Range("G3").Value = "0.7162"
$newSheet.Range("H3").Value = 6

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "005303"
$newSheet.Range("C4").Value = "嘉实医药健康股票A"
$newSheet.Range("D4").Value = "14.18"
$newSheet.Range("E4").Value = "92.14"
$newSheet.Range("F4").Value = "5.00"
$newSheet.Range("G4").Value = "0.7090"
$newSheet.Range("H4").Value = 7

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "001766"
$newSheet.Range("C5").Value = "上投摩根医疗健康股票"
$newSheet.Range("D5").Value = "10.35"
$newSheet.Range("E5").Value = "80.54"
$newSheet.Range("F5").Value = "6.02"
$newSheet.Range("G5").Value = "0.6231"
$newSheet.Range("H5").Value = 4

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "007553"
$newSheet.Range("C6").Value = "中信建投医改灵活配置混合C"
$newSheet.Range("D6").Value = "13.34"
$newSheet.Range("E6").Value = "94.92"
$newSheet.Range("F6").Value = "4.65"
$newSheet.Range("G6").Value = "0.6203"
$newSheet.Range("H6").Value = 8

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "010090"
$newSheet.Range("C7").Value = "中信建投医药健康混合A"
$newSheet.Range("D7").Value = "5.37"
$newSheet.Range("E7").Value = "94.87"
$newSheet.Range("F7").Value = "5.71"
$newSheet.Range("G7").Value = "0.3066"
$newSheet.Range("H7").Value = 7

$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "005108"
$newSheet.Range("C8").Value = "圆信永丰双利优选定期开放灵活配置混合"
$newSheet.Range("D8").Value = "1.89"
$newSheet.Range("E8").Value = "94.60"
$newSheet.Range("F8").Value = "9.36"
$newSheet.Range("G8").Value = "0.1769"
$newSheet.Range("H8").Value = 1

$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "005304"
$newSheet.Range("C9").Value = "嘉实医药健康股票C"
$newSheet.Range("D9").Value = "2.85"
$newSheet.Range("E9").Value = "92.14"
$newSheet.Range("F9").Value = "5.00"
$newSheet.Range("G9").Value = "0.1425"
$newSheet.Range("H9").Value = 7

$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "010091"
$newSheet.Range("C10").Value = "中信建投医药健康混合C"
$newSheet.Range("D10").Value = "2.25"
$newSheet.Range("E10").Value = "94.87"
$newSheet.Range("F10").Value = "5.71"
$newSheet.Range("G10").Value = "0.1285"
$newSheet.Range("H10").Value = 7

$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "008107"
$newSheet.Range("C11").Value = "华商医药医疗行业股票"
$newSheet.Range("D11").Value = "1.57"
$newSheet.Range("E11").Value = "91.35"
$newSheet.Range("F11").Value = "3.92"
$newSheet.Range("G11").Value = "0.0615"
$newSheet.Range("H11").Value = 9

$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "620001"
$newSheet.Range("C12").Value = "金元顺安宝石动力混合"
$newSheet.Range("D12").Value = "1.90"
$newSheet.Range("E12").Value = "56.14"
$newSheet.Range("F12").Value = "2.66"
$newSheet.Range("G12").Value = "0.0505"
$newSheet.Range("H12").Value = 8

$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "001965"
$newSheet.Range("C13").Value = "圆信永丰兴源灵活配置混合A"
$newSheet.Range("D13").Value = "0.76"
$newSheet.Range("E13").Value = "93.43"
$newSheet.Range("F13").Value = "6.18"
$newSheet.Range("G13").Value = "0.0470"
$newSheet.Range("H13").Value = 6

$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = "001966"
$newSheet.Range("C14").Value = "圆信永丰兴源灵活配置混合C"
$newSheet.Range("D14").Value = "0.25"
$newSheet.Range("E14").Value = "93.43"
$newSheet.Range("F14").Value = "6.18"
$newSheet.Range("G14").Value = "0.0154"
$newSheet.Range("H14").Value = 6

$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = "006274"
$newSheet.Range("C15").Value = "圆信永丰医药健康混合"
$newSheet.Range("D15").Value = "0.18"
$newSheet.Range("E15").Value = "93.66"
$newSheet.Range("F15").Value = "6.51"
$newSheet.Range("G15").Value = "0.0117"
$newSheet.Range("H15").Value = 4

$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = "004536"
$newSheet.Range("C16").Value = "嘉实中小企业量化活力灵活配置混合"
$newSheet.Range("D16").Value = "0.17"
$newSheet.Range("E16").Value = "90.06"
$newSheet.Range("F16").Value = "1.65"
$newSheet.Range("G16").Value = "0.0028"
$newSheet.Range("H16").Value = 5

# Now that "2022-Q1" is fully populated, move it into place: right before "总计"
$newSheet.Move($wb.Worksheets.Item("总计"))

# --- Update the "总计" (summary) sheet: insert the 2022-Q1 row at the top of the data ---
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()
$summary.Rows.Item(2).ClearFormats()
$summary.Range("A3").Copy($summary.Range("A2"))
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q1"
$summary.Range("C2").Value = 15
$summary.Range("D2").Value = 4.81

# Renumber the running index in column A for the rows that shifted down
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3

